# Tracker.xlsx — add the "Python" rows (13-15) covering the Interpreter/
# data-types/operators, compound data structures, and control-flow/functions
# sessions, matching the existing CORE TOPIC / DISCRIPTION / DATE table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone the formatting (date number format + colored fills/borders) of the
# last existing data row (row 12) onto the three new rows so no new styles
# are minted.
$ws.Range("A12:C12").Copy()
$ws.Range("A13:C13").PasteSpecial(-4122)
$ws.Range("A14:C14").PasteSpecial(-4122)
$ws.Range("A15:C15").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 13 - The Interpreter / environment / intro / data types / operators
$ws.Range("A13").Value = 43322
$ws.Range("B13").Value = "Python"
$ws.Range("C13").Value = "The Interpreter, And its Environment, Introducton to Python, Data Types, Operators"

# Row 14 - Compound data structures
$ws.Range("A14").Value = 43323
$ws.Range("B14").Value = "Python"
$ws.Range("C14").Value = " Lists, Tuples,Sets, Dictionaries, Compound data structures, and their Methods,"

# Row 15 - Control flow tools & functions
$ws.Range("A15").Value = 43324
$ws.Range("B15").Value = "Python"
$ws.Range("C15").Value = "Control Flow tools, functions"

$ws.Range("C15").Select()
